# Applies the "using gains for all" edit:
#  - inserts two new columns (F, G) before the existing "GFA - Sales" column,
#    shifting the old F:M columns to H:O
#  - adds header labels "M_TotalTax" (F1) and "M_CorpTax" (G1)
#  - fills in the new M_TotalTax / M_CorpTax data for rows 2-11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at F:G; this pushes old columns F:M to H:O
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New data values for M_TotalTax (F) and M_CorpTax (G), rows 2-11 (ordered by row)
$data = @(
    @(6308727034979.312, 399825921028.5854),
    @(16630145391623.02, 1639742485782.957),
    @(4450994137606.095, 601350231413.5104),
    @(4183547438952.192, 598849276038.3025),
    @(11223287075501.79, 872292028558.4308),
    @(1841737275230.086, 214321200777.9413),
    @(6192585801479.285, 516695167857.3162),
    @(14653861967257.56, 1232540278767.842),
    @(9623160693235.053, 876943418066.7275),
    @(5030701274022.499, 355596860701.1148)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
}
